# Authentication_TestData.xlsx - "Extent Report Configuration is done"
#
# The sheet holds several small "table blocks" (TestCaseName / Description /
# RunMode / RowID / expected-response columns). This edit:
#   - Renames the "Expected*" header labels to their shorter API field names
#     (ExpectedHttpCode -> getStatusCode, Expectedstatus_message ->
#     status_message, Expectedsuccess -> success, Expectedstatus_code ->
#     status_code) and moves them one column to the left (into column E)
#     for the first two blocks.
#   - Fixes a typo in a test-case id ("...With_Invalid_Resource_TestData"
#     -> "...With_invalid_Resource_TestData").
#   - Re-purposes the old "Invalid Resource" row (row 13) into a proper,
#     fully populated "Validate Expected Response with InValid Resource
#     Request" test-data row (404 / message / RowID 34), renumbering its
#     test case id from TC02 to TC03.
#   - Populates the previously-blank "getStatusCode" column (E) for the
#     CreateSession block (rows 17-20) with the expected http codes that
#     used to live in column F, and clears the now-unused trailing columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Block 1 (rows 1-3): TS01_TC01 ... Valid_Key -----------------------
$ws.Range("E2").Value = "getStatusCode"
$ws.Range("F2").Value = "status_message"
$ws.Range("G2").Value = "success"
$ws.Range("H2").Value = "status_code"

# ---- Block 2 (rows 6-8): TS01_TC02 ... invalid_Key ----------------------
$ws.Range("F7").Value = "getStatusCode"
$ws.Range("G7").Value = "status_message"
$ws.Range("H7").Value = "success"
$ws.Range("I7").Value = "status_code"

# ---- Block 3 (rows 11-13): TS01_TC03 ... invalid_Resource ---------------
$ws.Range("A11").Value = "TS01_TC03_Authencation_CreateTestToken_Get_With_invalid_Resource_TestData"

$ws.Range("E12").Value = "getStatusCode"
$ws.Range("F12").Value = "status_message"
$ws.Range("G12").Value = "status_code"
$ws.Range("H12").Clear()
$ws.Range("I12").Clear()

$ws.Range("A13").Value = "TS01_TC03_Authencation_CreateTestToken_Get_With_Invalid_Resource"
$ws.Range("B13").Value = "Validate Expected Response with  InValid Resource Request"
$ws.Range("E13").Value = "404"
$ws.Range("F13").Value = "The resource you requested could not be found."
$ws.Range("G13").Value = "34"
$ws.Range("H13").Clear()
$ws.Range("I13").Clear()

# ---- Block 4 (rows 16-20): Authencation_CreateSession_Get ---------------
$ws.Range("E17").Value = "getStatusCode"
$ws.Range("F17").Value = "status_message"
$ws.Range("G17").Value = "success"
$ws.Range("H17").Value = "status_code"
$ws.Range("I17").Clear()

$ws.Range("E18").Value = "200"
$ws.Range("F18").Clear()
$ws.Range("I18").Clear()

$ws.Range("E19").Value = "401"
$ws.Range("F19").Clear()
$ws.Range("I19").Clear()

$ws.Range("E20").Value = "404"
$ws.Range("F20").Clear()
$ws.Range("I20").Clear()

# ---- Misc: last active selection, as recorded by Excel on save ----------
$ws.Range("E28").Select()
